$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the recorded date string (A2) from 30.11.2022 -> 30.11.2023
$ws.Range("A2").Value = "30.11.2023"

# Move/restore the active selection to D8, matching the saved view state
$ws.Range("D8").Select() | Out-Null
